$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2452
$ws.Range("E2").Value = -113
$ws.Range("F2").Value = -113
$ws.Range("G2").Value = -245
$ws.Range("H2").Value = -242
$ws.Range("I2").Value = -236
$ws.Range("J2").Value = -6
$ws.Range("K2").Value = 5064
$ws.Range("L2").Value = 3297
$ws.Range("M2").Value = 1767
$ws.Range("N2").Value = 1715
$ws.Range("O2").Value = 51
$ws.Range("P2").Value = 84
$ws.Range("Q2").Value = -19
$ws.Range("R2").Value = -8
$ws.Range("S2").Value = 43
$ws.Range("T2").Value = 57
$ws.Range("U2").Value = -76
$ws.Range("V2").Value = 2683
$ws.Range("W2").Value = -4.61
$ws.Range("X2").Value = -9.859999999999999
$ws.Range("Y2").Value = -12.87
$ws.Range("Z2").Value = -4.69
$ws.Range("AA2").Value = 186.62
$ws.Range("AB2").Value = 2144.58
$ws.Range("AC2").Value = -14024
$ws.Range("AD2").Value = -2.45
$ws.Range("AE2").Value = 150535
$ws.Range("AF2").Value = 0.23
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 1680000

# Row 3
$ws.Range("D3").Value = 2373
$ws.Range("E3").Value = -105
$ws.Range("F3").Value = -105
$ws.Range("G3").Value = -234
$ws.Range("H3").Value = -235
$ws.Range("I3").Value = -235
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5195
$ws.Range("L3").Value = 3336
$ws.Range("M3").Value = 1859
$ws.Range("N3").Value = 1807
$ws.Range("O3").Value = 52
$ws.Range("P3").Value = 84
$ws.Range("Q3").Value = 20
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = -65
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 18
$ws.Range("V3").Value = 2623
$ws.Range("W3").Value = -4.43
$ws.Range("X3").Value = -9.890000000000001
$ws.Range("Y3").Value = -13.33
$ws.Range("Z3").Value = -4.58
$ws.Range("AA3").Value = 179.43
$ws.Range("AB3").Value = 1866.22
$ws.Range("AC3").Value = -13977
$ws.Range("AD3").Value = -4.39
$ws.Range("AE3").Value = 158604
$ws.Range("AF3").Value = 0.39
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 1680000

# Row 4
$ws.Range("D4").Value = 2008
$ws.Range("E4").Value = -125
$ws.Range("F4").Value = -125
$ws.Range("G4").Value = -173
$ws.Range("H4").Value = -175
$ws.Range("I4").Value = -169
$ws.Range("J4").Value = -6
$ws.Range("K4").Value = 4841
$ws.Range("L4").Value = 3179
$ws.Range("M4").Value = 1662
$ws.Range("N4").Value = 1616
$ws.Range("O4").Value = 46
$ws.Range("P4").Value = 84
$ws.Range("Q4").Value = 44
$ws.Range("R4").Value = 144
$ws.Range("S4").Value = -178
$ws.Range("T4").Value = 2
$ws.Range("U4").Value = 42
$ws.Range("V4").Value = 2472
$ws.Range("W4").Value = -6.21
$ws.Range("X4").Value = -8.73
$ws.Range("Y4").Value = -9.869999999999999
$ws.Range("Z4").Value = -3.49
$ws.Range("AA4").Value = 191.29
$ws.Range("AB4").Value = 1663.22
$ws.Range("AC4").Value = -10057
$ws.Range("AD4").Value = -2.61
$ws.Range("AE4").Value = 141856
$ws.Range("AF4").Value = 0.18
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 1680000

# Row 5
$ws.Range("D5").Value = 1825
$ws.Range("E5").Value = -181
$ws.Range("F5").Value = -181
$ws.Range("G5").Value = -151
$ws.Range("H5").Value = -130
$ws.Range("I5").Value = -120
$ws.Range("J5").Value = -10
$ws.Range("K5").Value = 4135
$ws.Range("L5").Value = 2646
$ws.Range("M5").Value = 1489
$ws.Range("N5").Value = 1489
$ws.Range("P5").Value = 84
$ws.Range("Q5").Value = -115
$ws.Range("R5").Value = 470
$ws.Range("S5").Value = -363
$ws.Range("T5").Value = 2
$ws.Range("U5").Value = -116
$ws.Range("V5").Value = 2075
$ws.Range("W5").Value = -9.949999999999999
$ws.Range("X5").Value = -7.15
$ws.Range("Y5").Value = -7.76
$ws.Range("Z5").Value = -2.9
$ws.Range("AA5").Value = 177.71
$ws.Range("AB5").Value = 1512.5
$ws.Range("AC5").Value = -7169
$ws.Range("AD5").Value = -2.69
$ws.Range("AE5").Value = 130677
$ws.Range("AF5").Value = 0.15
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 2.59
$ws.Range("AI5").Value = -4.73
$ws.Range("AJ5").Value = 1680000

# Row 6
$ws.Range("D6").Value = 1493
$ws.Range("E6").Value = -87
$ws.Range("F6").Value = -87
$ws.Range("G6").Value = -291
$ws.Range("H6").Value = -233
$ws.Range("I6").Value = -233
$ws.Range("K6").Value = 3640
$ws.Range("L6").Value = 2390
$ws.Range("M6").Value = 1250
$ws.Range("N6").Value = 1250
$ws.Range("P6").Value = 84
$ws.Range("Q6").Value = -122
$ws.Range("R6").Value = 286
$ws.Range("S6").Value = -173
$ws.Range("T6").Value = 5
$ws.Range("U6").Value = -127
$ws.Range("V6").Value = 1907
$ws.Range("W6").Value = -5.84
$ws.Range("X6").Value = -15.57
$ws.Range("Y6").Value = -16.98
$ws.Range("Z6").Value = -5.98
$ws.Range("AA6").Value = 191.25
$ws.Range("AB6").Value = 1298.52
$ws.Range("AC6").Value = -13840
$ws.Range("AD6").Value = -1.24
$ws.Range("AE6").Value = 109682
$ws.Range("AF6").Value = 0.16
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 1680000

# Row 5: remove O5 (deleted in target)
$ws.Range("O5").ClearContents()

# Row 6: remove AG6 and AH6 (deleted in target)
$ws.Range("AG6:AH6").ClearContents()

# Rows 7-9: clear all data cells (keep only A,B,C)
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()

Write-Output "edit applied"